$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting existing data from D:K to F:M
$ws.Columns("D:E").Insert()

# Copy number formatting (styles) from column F into new D:E columns,
# restricted to the rows that actually contain data (skipping the blank
# separator rows between the three statement blocks)
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new column D and E with the newest two quarters of data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1261700
$ws.Range("E8").Value = 1281100
$ws.Range("D9").Value = 746200
$ws.Range("E9").Value = 737500
$ws.Range("D10").Value = 515500
$ws.Range("E10").Value = 543600
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 17700
$ws.Range("E14").Value = -1700
$ws.Range("D15").Value = 177200
$ws.Range("E15").Value = 175100
$ws.Range("D17").Value = 1061800
$ws.Range("E17").Value = 1048200
$ws.Range("D18").Value = 199900
$ws.Range("E18").Value = 232900
$ws.Range("D20").Value = -4600
$ws.Range("E20").Value = 5200
$ws.Range("D21").Value = 372500
$ws.Range("E21").Value = 413200
$ws.Range("D22").Value = 29300
$ws.Range("E22").Value = 35200
$ws.Range("D23").Value = 166000
$ws.Range("E23").Value = 202900
$ws.Range("D24").Value = 33200
$ws.Range("E24").Value = 52100
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 132800
$ws.Range("E26").Value = 150800
$ws.Range("D27").Value = 132700
$ws.Range("E27").Value = 150800
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -300
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 4600
$ws.Range("E32").Value = -5200
$ws.Range("D33").Value = 132500
$ws.Range("E33").Value = 150800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 132500
$ws.Range("E35").Value = 150800
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 319300
$ws.Range("E41").Value = 244400
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 609500
$ws.Range("E43").Value = 625000
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 164100
$ws.Range("E45").Value = 153600
$ws.Range("D46").Value = 1092900
$ws.Range("E46").Value = 1023100
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 5169000
$ws.Range("E48").Value = 5069800
$ws.Range("D49").Value = 6160300
$ws.Range("E49").Value = 5882400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 205100
$ws.Range("E52").Value = 216900
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 12627300
$ws.Range("E54").Value = 12192100
$ws.Range("D57").Value = 360000
$ws.Range("E57").Value = 331500
$ws.Range("D58").Value = 20300
$ws.Range("E58").Value = 20800
$ws.Range("D59").Value = 480400
$ws.Range("E59").Value = 462900
$ws.Range("D60").Value = 860700
$ws.Range("E60").Value = 815200
$ws.Range("D61").Value = 4153500
$ws.Range("E61").Value = 3747200
$ws.Range("D62").Value = 1153000
$ws.Range("E62").Value = 1124500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 6172700
$ws.Range("E66").Value = 5692500
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 2264500
$ws.Range("E72").Value = 2174100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 6454600
$ws.Range("E76").Value = 6499600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 132500
$ws.Range("E81").Value = 150800
$ws.Range("D83").Value = 177200
$ws.Range("E83").Value = 175100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 373400
$ws.Range("E89").Value = 373000
$ws.Range("D91").Value = -172600
$ws.Range("E91").Value = -171800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -501400
$ws.Range("E94").Value = -185200
$ws.Range("D96").Value = -42100
$ws.Range("E96").Value = -36900
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 208100
$ws.Range("E100").Value = -86300
$ws.Range("D101").Value = -800
$ws.Range("E101").Value = 400
$ws.Range("D102").Value = 79400
$ws.Range("E102").Value = 101900

# Apply a few data corrections to previously existing quarters
$ws.Range("I91").Value = -114800
$ws.Range("J91").Value = -202600
$ws.Range("H94").Value = -177400
$ws.Range("H101").Value = 800
$ws.Range("H102").Value = 43900
